$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new value in D26 - "in progress"
$ws.Range("D26").Value = "in progress"

# Add new rows 29 and 30
$ws.Range("C29").Value = "make the create_word_xml work in unix/mac - it uses command line ot copy files. "
$ws.Range("C30").Value = "maybe move RESERVED to .reserved, which would be ignored in ls(), and might be cleaner? Or should I make it totally visiable.."

# Update the selection to match target
$ws.Range("C26").Select()
